$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet: updated weekly forecast figures (4wk low sales check) ---
$ws1.Range("D2").Value = 25
$ws1.Range("H2").Value = 12.96
$ws1.Range("L2").Value = 1.09
$ws1.Range("D3").Value = 25
$ws1.Range("H3").Value = 11.96
$ws1.Range("L3").Value = 1.04
$ws1.Range("D4").Value = 26
$ws1.Range("H4").Value = 10.54
$ws1.Range("L4").Value = 0.92
$ws1.Range("D5").Value = 27
$ws1.Range("H5").Value = 9.19
$ws1.Range("L5").Value = 0.85
$ws1.Range("D6").Value = 27
$ws1.Range("H6").Value = 8.19
$ws1.Range("L6").Value = 0.86
$ws1.Range("D7").Value = 27
$ws1.Range("H7").Value = 7.19
$ws1.Range("L7").Value = 1.01
$ws1.Range("D8").Value = 26
$ws1.Range("H8").Value = 6.42
$ws1.Range("L8").Value = 1.04
$ws1.Range("D9").Value = 27
$ws1.Range("H9").Value = 5.22
$ws1.Range("L9").Value = 0.92
$ws1.Range("D10").Value = 28
$ws1.Range("H10").Value = 4.07
$ws1.Range("L10").Value = 0.89
$ws1.Range("D11").Value = 28
$ws1.Range("H11").Value = 3.07
$ws1.Range("L11").Value = 0.92
$ws1.Range("D12").Value = 23
$ws1.Range("H12").Value = 2.52
$ws1.Range("I12").Value = "Low"
$ws1.Range("J12").Value = "Normal"
$ws1.Range("L12").Value = 0.87
$ws1.Range("D13").Value = 20
$ws1.Range("H13").Value = 1.75
$ws1.Range("I13").Value = "Low"
$ws1.Range("J13").Value = "Normal"
$ws1.Range("L13").Value = 0.86
$ws1.Range("D14").Value = 20
$ws1.Range("H14").Value = 0.75
$ws1.Range("I14").Value = "Low"
$ws1.Range("L14").Value = 0.95
$ws1.Range("D15").Value = 24
$ws1.Range("L15").Value = 1.13
$ws1.Range("D16").Value = 27
$ws1.Range("L16").Value = 1.09
$ws1.Range("L17").Value = 0.83

# --- Summary sheet: recalculated forecast totals ---
$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B9").Value = "408"
$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "210"
$ws2.Range("B11").NumberFormat = "@"
$ws2.Range("B11").Value = "103"
$ws2.Range("B12").NumberFormat = "@"
$ws2.Range("B12").Value = "28"
$ws2.Range("B14").NumberFormat = "@"
$ws2.Range("B14").Value = "20"

Write-Host "Applied 4wk low sales check updates"
